$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) values for rows 4-12, columns D, I, J, K, L, M, N, P, Q.
# Columns A,B,C,E,F,G,H,O,R are unchanged across these rows.
$data = @(
    @{ Row = 4;  D = 44536; I = "Primera"; J = 87; K = 22000; L = 22000; M = 22000; N = "`$/bandeja 18 kilos"; P = 1222; Q = 18 },
    @{ Row = 5;  D = 44536; I = "Segunda"; J = 80; K = 20000; L = 20000; M = 20000; N = "`$/bandeja 18 kilos"; P = 1111; Q = 18 },
    @{ Row = 6;  D = 44424; I = "Primera"; J = 75; K = 18000; L = 18000; M = 18000; N = "`$/caja 15 kilos";    P = 1200; Q = 15 },
    @{ Row = 7;  D = 44424; I = "Segunda"; J = 50; K = 12000; L = 12000; M = 12000; N = "`$/caja 15 kilos";    P = 800;  Q = 15 },
    @{ Row = 8;  D = 44235; I = "Primera"; J = 80; K = 14000; L = 14000; M = 14000; N = "`$/bandeja 18 kilos"; P = 778;  Q = 18 },
    @{ Row = 9;  D = 44235; I = "Segunda"; J = 70; K = 12000; L = 12000; M = 12000; N = "`$/bandeja 18 kilos"; P = 667;  Q = 18 },
    @{ Row = 10; D = 44235; I = "Tercera"; J = 60; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; P = 556;  Q = 18 },
    @{ Row = 11; D = 44242; I = "Primera"; J = 60; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 },
    @{ Row = 12; D = 44242; I = "Segunda"; J = 50; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; P = 556;  Q = 18 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D    # D: Fecha
    $ws.Cells.Item($r, 9).Value = $entry.I    # I: Calidad
    $ws.Cells.Item($r, 10).Value = $entry.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $entry.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $entry.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $entry.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $entry.N   # N: Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $entry.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $entry.Q   # Q: Kg o Unidades
}
